$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 34 values (revision of the "01-01-2021" quarter)
$row34 = @{
    'B' = 1974; 'C' = 229;  'D' = 3600; 'E' = 3237; 'F' = 362;  'G' = 4087;
    'H' = 1134; 'I' = 519;  'J' = 89;   'K' = 216;  'L' = 279;  'M' = 311;
    'N' = 625;  'O' = 231;  'P' = 674;  'Q' = 1055; 'R' = 1999; 'S' = 4161;
    'T' = 577;  'U' = 1714; 'V' = 1431; 'W' = 2214; 'X' = 3629; 'Y' = 2963;
    'Z' = 3315; 'AA' = 1882; 'AB' = 34660; 'AC' = 3210; 'AD' = 206; 'AE' = 38076
}

foreach ($col in $row34.Keys) {
    $ws.Range($col + "34").Value = $row34[$col]
}

# Add new row 35 for the "01-04-2021" quarter
# Force text entry so the date-like string isn't auto-converted to a date
# serial number (it must round-trip as a shared string, like the other
# cells in column A), then restore the default "Normal" style.
$ws.Range("A35").NumberFormat = "@"
$ws.Range("A35").Value = "01-04-2021"
$ws.Range("A35").Style = "Normal"

$row35 = @{
    'B' = 870;  'C' = 272;  'D' = 3775; 'E' = 3368; 'F' = 406;  'G' = 4145;
    'H' = 1090; 'I' = 600;  'J' = 89;   'K' = 220;  'L' = 287;  'M' = 270;
    'N' = 615;  'O' = 248;  'P' = 726;  'Q' = 1074; 'R' = 2160; 'S' = 4066;
    'T' = 611;  'U' = 1565; 'V' = 1426; 'W' = 2223; 'X' = 3902; 'Y' = 2902;
    'Z' = 5176; 'AA' = 1904; 'AB' = 36066; 'AC' = 3454; 'AD' = 165; 'AE' = 39677
}

foreach ($col in $row35.Keys) {
    $ws.Range($col + "35").Value = $row35[$col]
}
